# Template-ify the order sheet and add a blank second sheet.
# (EasyExcel 读写 excel 文件: turn the sample-filled order sheet into a
# blank fill-in template, and add an empty Sheet2.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the three sample order rows (A3:A5 numbering + B3:E5 product/qty/
# price/date) back to blank template rows, same as rows 6-18 below them.
# Use Clear() on column A so the now-default-styled cells drop out of the
# sheet entirely (matching the rest of the unfilled rows), and
# ClearContents() on B:E so their per-column formatting (text/date/number)
# is retained for future data entry.
$ws.Range("A3:A5").Clear()
$ws.Range("B3:E5").ClearContents()

# Reset the active selection from E6 to A6.
[void]$ws.Range("A6").Select()

# Add a second, empty worksheet right after Sheet1, then restore Sheet1
# as the active/selected tab.
[void]$wb.Worksheets.Add($null, $ws)
[void]$ws.Activate()
